# Update the "想去人数" (want-to-go count) values on both the "展览"
# and "全部类型" sheets, which carry duplicate copies of the exhibition
# data rows.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 78
    $ws.Range("F3").Value = 7
}
